# Renumber the Word-generated "_Toc" heading bookmarks (Table of Contents
# anchors) to match a refreshed TOC field numbering sequence.
#
# Word COM does not allow directly writing Bookmark.Name (it is effectively
# read-only through automation for these hidden "_Toc" bookmarks), so each
# bookmark is recreated in-place: remember its Range, delete the bookmark,
# then add a new bookmark with the new name over the exact same Range. This
# preserves the bookmark's w:id ordering/position while only changing
# w:name, matching how Word itself updates these anchors when a TOC is
# rebuilt.

$d = $word.ActiveDocument

$renames = @(
    @{ Old = "_Toc5728518"; New = "_Toc5791597" },
    @{ Old = "_Toc5728519"; New = "_Toc5791598" },
    @{ Old = "_Toc5728520"; New = "_Toc5791599" },
    @{ Old = "_Toc5728521"; New = "_Toc5791600" },
    @{ Old = "_Toc5728522"; New = "_Toc5791601" },
    @{ Old = "_Toc5728523"; New = "_Toc5791602" },
    @{ Old = "_Toc5728524"; New = "_Toc5791603" },
    @{ Old = "_Toc5728525"; New = "_Toc5791604" },
    @{ Old = "_Toc5728526"; New = "_Toc5791605" },
    @{ Old = "_Toc5728527"; New = "_Toc5791606" },
    @{ Old = "_Toc5728528"; New = "_Toc5791607" },
    @{ Old = "_Toc5728529"; New = "_Toc5791608" },
    @{ Old = "_Toc5728530"; New = "_Toc5791609" }
)

foreach ($pair in $renames) {
    $bm = $d.Bookmarks.Item($pair.Old)
    $r = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($pair.New, $r) | Out-Null
}
